# Update the workbook "Avverkningsanmälningar" sheet:
#  1. Change every "Förändrad" (column C) date value from 45186 to 45188
#     for all data rows (rows 2 through 463).
#  2. Give row 463 an explicit row height (15pt, custom height) like the
#     other data rows already have.
#  3. Append a new data row (row 464) with a new logging notification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C (Förändrad) for rows 2..463 from 45186 -> 45188
for ($r = 2; $r -le 463; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45188
}

# 2. Row 463 gains an explicit custom row height (matches ht="15" customHeight="1")
$ws.Rows.Item(463).RowHeight = 15

# 3. Append new row 464 with the new record
$row = 464

$ws.Cells.Item($row, 1).Value2 = "A 43858-2023"      # A - Beteckning

$ws.Cells.Item($row, 2).Value2 = 45187               # B - Datum
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 3).Value2 = 45188               # C - Förändrad
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 4).Value2 = "JÖNKÖPINGS LÄN"    # D - Län
$ws.Cells.Item($row, 5).Value2 = "VÄRNAMO"           # E - Kommun
# F (Markägare) intentionally left blank

$ws.Cells.Item($row, 7).Value2 = 3.3                 # G - Area (ha)
$ws.Cells.Item($row, 8).Value2 = 0                   # H - Fridlysta
$ws.Cells.Item($row, 9).Value2 = 0                   # I - Signalarter
$ws.Cells.Item($row, 10).Value2 = 0                  # J - NT
$ws.Cells.Item($row, 11).Value2 = 0                  # K - VU
$ws.Cells.Item($row, 12).Value2 = 0                  # L - EN
$ws.Cells.Item($row, 13).Value2 = 0                  # M - CR
$ws.Cells.Item($row, 14).Value2 = 0                  # N - RE
$ws.Cells.Item($row, 15).Value2 = 0                  # O - Rödlistade
$ws.Cells.Item($row, 16).Value2 = 0                  # P - Hotade
$ws.Cells.Item($row, 17).Value2 = 0                  # Q - Alla arter

# R - Artnamn: empty cell with wrap-text style (same as the other rows)
$ws.Cells.Item($row, 18).WrapText = $true
